$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 84 into a new row 85, preserving original cell types/values.
$ws.Range("A84:H84").Copy($ws.Range("A85:H85"))

# New row 85 keeps the same data as the old row 84, except the issue_type
# (column E) is updated from "WRI" to "EXP".
$ws.Range("E85").Value = "EXP"

# Original row 84 keeps its values, but politeness_score (column B) becomes
# a real number instead of a text string.
$ws.Range("B84").Value = 3
